$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.04941422446394483, 0.01396578115299896, 0.4181095735075786, 0.8541194379906329, 0.7137494338910955, 0.7440699538300066, 1.815162008128652)
    3 = @(0.0438473038125835, 0.01430705402919941, 0.3644670215709738, 0.8381558220722667, 0.6968537578041207, 0.7439463256324643, 1.597421177966396)
    4 = @(0.04044776810401629, 0.01453142143158992, 0.3316615277020674, 0.8293163305418432, 0.6874178380439844, 0.7445721242476537, 1.46373772809585)
    5 = @(0.03906697724346486, 0.0146265512858097, 0.3183231727715139, 0.8259536445128788, 0.6838055197740545, 0.7450023449477072, 1.409262595861151)
    6 = @(0.03883796906031023, 0.01464257005850378, 0.316110087367079, 0.8254096645835887, 0.6832196760267522, 0.7450843282223616, 1.400217143781617)
    7 = @(0.04042912802262322, 0.01453268945055086, 0.331481523621008, 0.8292700138362363, 0.6873681820376447, 0.7445772186314485, 1.463003050226575)
    8 = @(0.0474908182388134, 0.01408035715871137, 0.3995846383952539, 0.8484139667033475, 0.7077272898518601, 0.7438809324301161, 1.740082355985237)
    9 = @(0.06149248972270982, 0.0133123349444908, 0.5343119085064529, 0.8937054838520169, 0.7552330979233659, 0.7481408311370501, 2.283576567981072)
    10 = @(0.07188307377042236, 0.01282258809499481, 0.6342122584121199, 0.9318706771086056, 0.7949537421017396, 0.7547816025657994, 2.683099603312371)
    11 = @(0.07663472363594792, 0.01261642600701407, 0.6799026171229912, 0.950330828239359, 0.8141127874029337, 0.758582955027407, 2.86493305556354)
    12 = @(0.07843779455116362, 0.01254078891767563, 0.6972431061917348, 0.9574820664665964, 0.8215279911039204, 0.7601361028348208, 2.933803295394114)
    13 = @(0.07804930306535596, 0.01255696987882615, 0.6935067545089737, 0.95593472374307, 0.8199238204664425, 0.7597965266952542, 2.918970221198492)
    14 = @(0.07678298791140037, 0.01261015431531476, 0.6813284381240834, 0.950915926507264, 0.8147196128939527, 0.7587084472908145, 2.870598766937519)
    15 = @(0.0760078225384575, 0.01264304929070548, 0.6738739855137936, 0.9478627901375631, 0.8115528340603362, 0.758056811679154, 2.840971695168321)
    16 = @(0.07157305571345773, 0.01283639931576275, 0.6312314991648833, 0.9306866084641996, 0.7937238719233619, 0.7545490016616441, 2.671218233804325)
    17 = @(0.06885895203296855, 0.01295930043042048, 0.6051368309624365, 0.9204329421700521, 0.7830679169987036, 0.7525979757915877, 2.567103405680371)
    18 = @(0.06730020850478979, 0.01303155227081909, 0.5901507461954054, 0.9146385421048535, 0.7770412797225674, 0.7515491686880296, 2.507227556936925)
    19 = @(0.0667728418913498, 0.01305628261321701, 0.5850805456990429, 0.9126943040794941, 0.7750182469764866, 0.7512066188120343, 2.486956006867388)
    20 = @(0.06914763010942693, 0.01294605540508798, 0.6079122615375496, 0.9215137575612147, 0.7841916407862186, 0.7527980618409345, 2.578185745468375)
    21 = @(0.07715483315915606, 0.01259446641569184, 0.6849044323611508, 0.9523856829389388, 0.8162438431045587, 0.7590249464373073, 2.884806249797293)
    22 = @(0.08240977804257454, 0.01237887869179133, 0.7354494868061465, 0.9735007378138363, 0.8381265142589314, 0.7637577385979739, 3.085281863944545)
    23 = @(0.07960307552674806, 0.01249262892457281, 0.7084508395518299, 0.9621444330447275, 0.8263606768373677, 0.7611706012031334, 2.978276399807271)
    24 = @(0.0690171136829747, 0.01295203851863747, 0.606657439512901, 0.9210248078835122, 0.7836832952427528, 0.7527073760924452, 2.573175478438884)
    25 = @(0.0576872010046543, 0.01350716756799564, 0.4977182546596026, 0.8806062676985107, 0.741549675448141, 0.7463777939548208, 2.136518613275427)
}

$cols = @("C","D","E","F","G","H","K")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = [double]$vals[$i]
    }
}
